$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to remain plain text so values like "0.3730" or
# "28.561.66" are not reinterpreted as numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.561.66'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").Value = '1.795.02'
$ws.Range("E3").Value = '  -0.89%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.62%  '
$ws.Range("D5").Value = '327.73'
$ws.Range("E5").Value = '  -3.27%  '
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("D7").Value = '0.4391'
$ws.Range("E7").Value = '  -3.68%  '
$ws.Range("D8").Value = '0.3730'
$ws.Range("E8").Value = '  +5.32%  '
$ws.Range("D9").Value = '45.66'
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").Value = '0.07572'
$ws.Range("E10").Value = '  -0.52%  '
$ws.Range("D11").Value = '1.131'
$ws.Range("E11").Value = '  -2.20%  '
$ws.Range("D12").Value = '22.56'
$ws.Range("E12").Value = '  -1.24%  '
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("D14").Value = '6.193'
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("D15").Value = '7.468'
$ws.Range("E15").Value = '  +2.67%  '
$ws.Range("D16").Value = '1.800.66'
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("D17").Value = '0.00001085'
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D18").Value = '0.06712'
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("D19").Value = '80.52'
$ws.Range("E19").Value = '  -1.91%  '
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("D21").Value = '17.46'
$ws.Range("E21").Value = '  +0.83%  '
$ws.Range("D22").Value = '6.216'
$ws.Range("E22").Value = '  -3.06%  '
$ws.Range("D23").Value = '28.576.86'
$ws.Range("E23").Value = '  +1.07%  '
$ws.Range("D24").Value = '11.66'
$ws.Range("E24").Value = '  -2.98%  '
$ws.Range("D25").Value = '2.433'
$ws.Range("E25").Value = '  +1.08%  '
$ws.Range("D26").Value = '20.41'
$ws.Range("E26").Value = '  -2.23%  '
$ws.Range("D27").Value = '153.19'
$ws.Range("E27").Value = '  -1.54%  '
$ws.Range("D28").Value = '2.331'
$ws.Range("E28").Value = '  -4.45%  '
$ws.Range("D29").Value = '2.009.83'
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("D30").Value = '1.311'
$ws.Range("E30").Value = '  +0.49%  '
$ws.Range("D31").Value = '130.35'
$ws.Range("E31").Value = '  -2.71%  '
$ws.Range("D32").Value = '3.978'
$ws.Range("E32").Value = '  -2.40%  '
$ws.Range("D33").Value = '5.763'
$ws.Range("E33").Value = '  -3.65%  '
$ws.Range("D34").Value = '0.09206'
$ws.Range("E34").Value = '  -3.06%  '
$ws.Range("D35").Value = '0.2235'
$ws.Range("E35").Value = '  +2.38%  '
$ws.Range("D36").Value = '12.06'
$ws.Range("E36").Value = '  -1.74%  '
$ws.Range("D37").Value = '0.06266'
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").Value = '0.02305'
$ws.Range("E38").Value = '  -3.75%  '
$ws.Range("D39").Value = '5.168'
$ws.Range("E39").Value = '  -1.02%  '
$ws.Range("D40").Value = '0.6564'
$ws.Range("E40").Value = '  -3.39%  '
$ws.Range("D41").Value = '1.192'
$ws.Range("E41").Value = '  -2.36%  '
$ws.Range("B42").Value = 'WEMIXTOKEN'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").Value = '1.426'
$ws.Range("E42").Value = '  -4.44%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '7.993'
$ws.Range("E43").Value = '  -2.81%  '
$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  +0.42%  '
$ws.Range("D45").Value = '13.84'
$ws.Range("E45").Value = '  -1.49%  '
$ws.Range("D46").Value = '0.6048'
$ws.Range("E46").Value = '  -1.70%  '
$ws.Range("D47").Value = '3.792'
$ws.Range("E47").Value = '  -2.03%  '
$ws.Range("D48").Value = '127.26'
$ws.Range("E48").Value = '  -2.37%  '
$ws.Range("D49").Value = '2.005'
$ws.Range("E49").Value = '  -2.18%  '
$ws.Range("D50").Value = '0.07003'
$ws.Range("E50").Value = '  -1.57%  '
$ws.Range("D51").Value = '77.46'
$ws.Range("E51").Value = '  -1.18%  '
